$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Indent the pseudo-code lines under "begin if" to visually tabulate them
$ws.Range("B23").Value = "   injured = true"
$ws.Range("B24").Value = "   then generate randomHealthIncrease(1-10)"
$ws.Range("B25").Value = "    playerHealthyTest(randomHealthIncrease int, playerId int)"

# Move the trailing comment from column F to column G on row 25
$comment = $ws.Range("F25").Value()
$ws.Range("F25").ClearContents()
$ws.Range("G25").Value = $comment

# Update selection to match the final state
$ws.Range("E28").Select()
